# Add new worksheet "Abu-PE" positioned after "PFOA-Tia" (before "PE_abu")
$wb = $excel.ActiveWorkbook
$afterSheet = $wb.Worksheets.Item("PFOA-Tia")
$ws = $wb.Worksheets.Add($null, $afterSheet)
$ws.Name = "Abu-PE"

# Header row
$ws.Range("A1").Value = 'Replicates'
$ws.Range("B1").Value = 'Genotype'
$ws.Range("C1").Value = 'Treatment'
$ws.Range("D1").Value = 'Age_maturity'
$ws.Range("E1").Value = 'Size_maturity'
$ws.Range("F1").Value = 'Day_1brood'
$ws.Range("G1").Value = 'Number_1brood'
$ws.Range("H1").Value = 'Day_2brood'
$ws.Range("I1").Value = 'Number_2brood'
$ws.Range("J1").Value = 'Interval_btwnbrds'
$ws.Range("K1").Value = 'Fecundity'
$ws.Range("J1").Interior.Color = 65535

# Data rows
# row 2
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 'LRV_0_1'
$ws.Range("C2").Value = 'Control'
$ws.Range("D2").Value = 7
$ws.Range("E2").Value = 3385.58
$ws.Range("F2").Value = 10
$ws.Range("G2").Value = 29
$ws.Range("H2").Value = 13
$ws.Range("I2").Value = 33
$ws.Range("J2").Formula = "=H2-F2"
$ws.Range("K2").Value = 62
# row 3
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 'LRV_0_1'
$ws.Range("C3").Value = 'Control'
$ws.Range("D3").Value = 7
$ws.Range("E3").Value = 3641.01
$ws.Range("F3").Value = 9
$ws.Range("G3").Value = 30
$ws.Range("H3").Value = 12
$ws.Range("I3").Value = 35
$ws.Range("J3").Formula = "=H3-F3"
$ws.Range("K3").Value = 65
# row 4
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 'LRV_0_1'
$ws.Range("C4").Value = 'Control'
$ws.Range("D4").Value = 7
$ws.Range("E4").Value = 3496.68
$ws.Range("F4").Value = 10
$ws.Range("G4").Value = 18
$ws.Range("H4").Value = 13
$ws.Range("I4").Value = 28
$ws.Range("J4").Value = 3
$ws.Range("K4").Value = 46
# row 5
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 'LRV_0_1'
$ws.Range("C5").Value = 'Control'
$ws.Range("D5").Value = 7
$ws.Range("E5").Value = 3695.36
$ws.Range("F5").Value = 9
$ws.Range("G5").Value = 27
$ws.Range("H5").Value = 12
$ws.Range("I5").Value = 31
$ws.Range("J5").Value = 3
$ws.Range("K5").Value = 58
# row 6
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = 'LRV_0_1'
$ws.Range("C6").Value = 'PE'
$ws.Range("D6").Value = 11
$ws.Range("E6").Value = 3399.68
$ws.Range("F6").Value = 13
$ws.Range("G6").Value = 6
$ws.Range("A6:K6").Interior.Color = 65535
# row 7
$ws.Range("A7").Value = 2
$ws.Range("B7").Value = 'LRV_0_1'
$ws.Range("C7").Value = 'PE'
$ws.Range("D7").Value = 11
$ws.Range("E7").Value = 3426.91
$ws.Range("F7").Value = 13
$ws.Range("G7").Value = 10
$ws.Range("H7").Value = 16
$ws.Range("I7").Value = 25
$ws.Range("J7").Value = 3
$ws.Range("K7").Value = 35
# row 8
$ws.Range("A8").Value = 3
$ws.Range("B8").Value = 'LRV_0_1'
$ws.Range("C8").Value = 'PE'
$ws.Range("D8").Value = 9
$ws.Range("E8").Value = 3440.28
$ws.Range("F8").Value = 11
$ws.Range("G8").Value = 17
$ws.Range("H8").Value = 13
$ws.Range("I8").Value = 22
$ws.Range("J8").Value = 2
$ws.Range("K8").Value = 39
# row 9
$ws.Range("A9").Value = 4
$ws.Range("B9").Value = 'LRV_0_1'
$ws.Range("C9").Value = 'PE'
$ws.Range("D9").Value = 11
$ws.Range("E9").Value = 3312.46
$ws.Range("F9").Value = 13
$ws.Range("G9").Value = 20
$ws.Range("H9").Value = 17
$ws.Range("I9").Value = 20
$ws.Range("J9").Value = 4
$ws.Range("K9").Value = 40
# row 10
$ws.Range("A10").Value = 1
$ws.Range("B10").Value = 'LR2_36_1'
$ws.Range("C10").Value = 'Control'
$ws.Range("D10").Value = 8
$ws.Range("E10").Value = 3696.1
$ws.Range("F10").Value = 11
$ws.Range("G10").Value = 28
$ws.Range("H10").Value = 14
$ws.Range("I10").Value = 30
$ws.Range("J10").Value = 3
$ws.Range("K10").Value = 58
# row 11
$ws.Range("A11").Value = 2
$ws.Range("B11").Value = 'LR2_36_1'
$ws.Range("C11").Value = 'Control'
$ws.Range("D11").Value = 8
$ws.Range("E11").Value = 3213.75
$ws.Range("F11").Value = 10
$ws.Range("G11").Value = 27
$ws.Range("H11").Value = 13
$ws.Range("I11").Value = 28
$ws.Range("J11").Value = 3
$ws.Range("K11").Value = 55
# row 12
$ws.Range("A12").Value = 3
$ws.Range("B12").Value = 'LR2_36_1'
$ws.Range("C12").Value = 'Control'
$ws.Range("D12").Value = 9
$ws.Range("A12:K12").Interior.Color = 65535
# row 13
$ws.Range("A13").Value = 4
$ws.Range("B13").Value = 'LR2_36_1'
$ws.Range("C13").Value = 'Control'
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = 3585.97
$ws.Range("F13").Value = 10
$ws.Range("G13").Value = 29
$ws.Range("H13").Value = 13
$ws.Range("I13").Value = 32
$ws.Range("J13").Value = 3
$ws.Range("K13").Value = 61
# row 14
$ws.Range("A14").Value = 1
$ws.Range("B14").Value = 'LR2_36_1'
$ws.Range("C14").Value = 'PE'
$ws.Range("D14").Value = 12
$ws.Range("A14:K14").Interior.Color = 65535
# row 15
$ws.Range("A15").Value = 2
$ws.Range("B15").Value = 'LR2_36_1'
$ws.Range("C15").Value = 'PE'
$ws.Range("D15").Value = 13
$ws.Range("E15").Value = 3571.93
$ws.Range("F15").Value = 16
$ws.Range("G15").Value = 16
$ws.Range("H15").Value = 19
$ws.Range("I15").Value = 19
$ws.Range("J15").Value = 3
$ws.Range("K15").Value = 35
# row 16
$ws.Range("A16").Value = 3
$ws.Range("B16").Value = 'LR2_36_1'
$ws.Range("C16").Value = 'PE'
$ws.Range("D16").Value = 8
$ws.Range("E16").Value = 2417.08
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 20
$ws.Range("H16").Value = 13
$ws.Range("I16").Value = 21
$ws.Range("J16").Value = 3
$ws.Range("K16").Value = 41
# row 17
$ws.Range("A17").Value = 4
$ws.Range("B17").Value = 'LR2_36_1'
$ws.Range("C17").Value = 'PE'
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = 3459.21
$ws.Range("F17").Value = 10
$ws.Range("G17").Value = 22
$ws.Range("H17").Value = 13
$ws.Range("I17").Value = 27
$ws.Range("J17").Value = 3
$ws.Range("K17").Value = 49

# Leave the cursor where the original author left it
$ws.Range("O18").Select()

# PFOA-Tia: re-enter the Interval_brood formulas for rows 6-14 as one range
# assignment so Excel stores them as a shared formula group (matches the
# upstream resave behaviour captured in the diff).
$wsPFOA = $wb.Worksheets.Item("PFOA-Tia")
$wsPFOA.Range("J6:J14").Formula = "=(F6+I6)"
